$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title
$ws.Range("B1").Value2 = "Nutrition Information Data Analysis Tool"

# Period selected (H2)
$ws.Range("H2").Value2 = 9

# Column B width (engine snaps to 1/7 increments; 40.7 is the closest achievable to 41.375)
$ws.Columns.Item(2).ColumnWidth = 40.7

# Activity rows B5:G39 (WBS relabeling + new estimate/actual data)
$data = New-Object 'object[,]' 35,6
$data[0,0] = "1.0 Initiating"
$data[0,1] = 1
$data[0,2] = 2
$data[0,3] = 1
$data[0,4] = 1
$data[0,5] = 1
$data[1,0] = "1.1 Develop Project Overview"
$data[1,1] = 1
$data[1,2] = 1
$data[1,3] = 1
$data[1,4] = 1
$data[1,5] = 0.5
$data[2,0] = "1.2 Develop Communication Plan"
$data[2,1] = 1
$data[2,2] = 1
$data[2,3] = 1
$data[2,4] = 1
$data[2,5] = 1
$data[3,0] = "2.0 Planning"
$data[3,1] = 2
$data[3,2] = 15
$data[3,3] = 1
$data[3,4] = 16
$data[3,5] = 0.2
$data[4,0] = "2.1 Develop System Vision Document"
$data[4,1] = 2
$data[4,2] = 2
$data[4,3] = $null
$data[4,4] = $null
$data[4,5] = 0
$data[5,0] = "2.2 Establish WBS, Activity Estimate and Gantt"
$data[5,1] = 2
$data[5,2] = 3
$data[5,3] = 2
$data[5,4] = 7
$data[5,5] = 1
$data[6,0] = "2.3 Define Project Requirements"
$data[6,1] = 5
$data[6,2] = 3
$data[6,3] = $null
$data[6,4] = $null
$data[6,5] = 0
$data[7,0] = "2.4 Define System Architecture"
$data[7,1] = 14
$data[7,2] = 2
$data[7,3] = $null
$data[7,4] = $null
$data[7,5] = 0
$data[8,0] = "2.5 Design UI Wireframes and Mockups"
$data[8,1] = 9
$data[8,2] = 3
$data[8,3] = $null
$data[8,4] = $null
$data[8,5] = 0
$data[9,0] = "3.0 Execution"
$data[9,1] = 17
$data[9,2] = 21
$data[9,3] = $null
$data[9,4] = $null
$data[9,5] = 0
$data[10,0] = "3.1 implement Food Search Feature"
$data[10,1] = 17
$data[10,2] = 7
$data[10,3] = $null
$data[10,4] = $null
$data[10,5] = 0
$data[11,0] = "3.2 Implement Nutrition Breakdown Feature"
$data[11,1] = 17
$data[11,2] = 7
$data[11,3] = $null
$data[11,4] = $null
$data[11,5] = 0
$data[12,0] = "3.3 implement Nutrition Range Filter Feature"
$data[12,1] = 24
$data[12,2] = 7
$data[12,3] = $null
$data[12,4] = $null
$data[12,5] = 0
$data[13,0] = "3.4 Implement Nutrition Level Feature"
$data[13,1] = 24
$data[13,2] = 7
$data[13,3] = $null
$data[13,4] = $null
$data[13,5] = 0
$data[14,0] = "3.5 Implement (TBD)"
$data[14,1] = 24
$data[14,2] = 7
$data[14,3] = $null
$data[14,4] = $null
$data[14,5] = 0
$data[15,0] = "3.6 Testing"
$data[15,1] = 31
$data[15,2] = 7
$data[15,3] = $null
$data[15,4] = $null
$data[15,5] = 0
$data[16,0] = "3.6.1 Unit Testing"
$data[16,1] = 31
$data[16,2] = 3
$data[16,3] = $null
$data[16,4] = $null
$data[16,5] = 0
$data[17,0] = "3.6.2 User Acceptance Testing"
$data[17,1] = 34
$data[17,2] = 4
$data[17,3] = $null
$data[17,4] = $null
$data[17,5] = 0
$data[18,0] = "4.0 Controlling"
$data[18,1] = 1
$data[18,2] = 44
$data[18,3] = 1
$data[18,4] = 44
$data[18,5] = 0.2
$data[19,0] = "4.1 Project Monitoring"
$data[19,1] = 1
$data[19,2] = 44
$data[19,3] = 1
$data[19,4] = 44
$data[19,5] = 0.2
$data[20,0] = "4.1.1 Weekly Team Check-In"
$data[20,1] = 1
$data[20,2] = 44
$data[20,3] = 1
$data[20,4] = 44
$data[20,5] = 0.2
$data[21,0] = "4.1.2 Document Meeting Outcome"
$data[21,1] = 1
$data[21,2] = 44
$data[21,3] = 1
$data[21,4] = 44
$data[21,5] = 0.2
$data[22,0] = "4.2 Performance Measurement"
$data[22,1] = 17
$data[22,2] = 28
$data[22,3] = $null
$data[22,4] = $null
$data[22,5] = 0
$data[23,0] = "4.2.1 Task Tracking"
$data[23,1] = 17
$data[23,2] = 28
$data[23,3] = $null
$data[23,4] = $null
$data[23,5] = 0
$data[24,0] = "4.2.2 Quality Checking"
$data[24,1] = 17
$data[24,2] = 28
$data[24,3] = $null
$data[24,4] = $null
$data[24,5] = 0
$data[25,0] = "4.3 Change Management"
$data[25,1] = 17
$data[25,2] = 28
$data[25,3] = $null
$data[25,4] = $null
$data[25,5] = 0
$data[26,0] = "4.3.1 Change Discussion"
$data[26,1] = 17
$data[26,2] = 28
$data[26,3] = $null
$data[26,4] = $null
$data[26,5] = 0
$data[27,0] = "4.3.2 Change Implementation"
$data[27,1] = 17
$data[27,2] = 28
$data[27,3] = $null
$data[27,4] = $null
$data[27,5] = 0
$data[28,0] = "5.0 Closing"
$data[28,1] = 38
$data[28,2] = 7
$data[28,3] = $null
$data[28,4] = $null
$data[28,5] = 0
$data[29,0] = "5.1 Documentation"
$data[29,1] = 38
$data[29,2] = 6
$data[29,3] = $null
$data[29,4] = $null
$data[29,5] = 0
$data[30,0] = "5.2 Final Performance Review"
$data[30,1] = 44
$data[30,2] = 1
$data[30,3] = $null
$data[30,4] = $null
$data[30,5] = 0
$data[31,0] = "5.3 Project Closure"
$data[31,1] = 44
$data[31,2] = 1
$data[31,3] = $null
$data[31,4] = $null
$data[31,5] = 0
$data[32,0] = "Activity 33"
$data[32,1] = $null
$data[32,2] = $null
$data[32,3] = $null
$data[32,4] = $null
$data[32,5] = 0
$data[33,0] = "Activity 34"
$data[33,1] = $null
$data[33,2] = $null
$data[33,3] = $null
$data[33,4] = $null
$data[33,5] = 0
$data[34,0] = "Activity 35"
$data[34,1] = $null
$data[34,2] = $null
$data[34,3] = $null
$data[34,4] = $null
$data[34,5] = 0
$ws.Range("B5:G39").Value2 = $data

# Selection (matches the new view state)
$ws.Range("M14").Select()
